$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5453.125
$ws.Range("I40").Value = 3806.875
$ws.Range("J40").Value = 7099.375
$ws.Range("K40").Value = 3806.875
$ws.Range("L40").Value = 7099.375
$ws.Range("M40").Value = -3631.875
$ws.Range("N40").Value = -7449.375

$ws.Range("H51").Value = 4412.609
$ws.Range("I51").Value = 3233.1667
$ws.Range("K51").Value = 3233.1667
$ws.Range("M51").Value = -2749.1667

$ws.Range("H70").Value = 4819.7827
$ws.Range("I70").Value = 4676.0557
$ws.Range("J70").Value = 5337.2
$ws.Range("K70").Value = 14028.1671
$ws.Range("L70").Value = 16011.6
$ws.Range("M70").Value = -13758.1671
$ws.Range("N70").Value = -16551.6

$ws.Range("H73").Value = 4819.7827
$ws.Range("I73").Value = 4676.0557
$ws.Range("J73").Value = 5337.2
$ws.Range("K73").Value = 14028.1671
$ws.Range("L73").Value = 16011.6
$ws.Range("M73").Value = -13092.1671
$ws.Range("N73").Value = -17883.6

$ws.Range("H132").Value = 3323.879
$ws.Range("I132").Value = 3271.5625
$ws.Range("K132").Value = 9814.6875
$ws.Range("M132").Value = -7284.6875

$ws.Range("H138").Value = 83338350
$ws.Range("J138").Value = 83338350
$ws.Range("L138").Value = 250015050
$ws.Range("N138").Value = -250025330

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H95").Value = 67069
$ws.Range("J95").Value = 67069
$ws.Range("L95").Value = 67069
$ws.Range("N95").Value = -72561

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1404.9048
$ws.Range("I64").Value = 1561.2
$ws.Range("J64").Value = 1262.8182
$ws.Range("K64").Value = 1561.2
$ws.Range("L64").Value = 1262.8182
$ws.Range("M64").Value = -1336.2
$ws.Range("N64").Value = -1712.8182

$ws.Range("H67").Value = 1404.9048
$ws.Range("I67").Value = 1561.2
$ws.Range("J67").Value = 1262.8182
$ws.Range("K67").Value = 1561.2
$ws.Range("L67").Value = 1262.8182
$ws.Range("M67").Value = -781.2
$ws.Range("N67").Value = -2822.8182

$ws.Range("H117").Value = 38499.5
$ws.Range("J117").Value = 38499.5
$ws.Range("L117").Value = 38499.5
$ws.Range("N117").Value = -47677.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2709
$ws.Range("I58").Value = 2367.85
$ws.Range("J58").Value = 4983.3335
$ws.Range("K58").Value = 2367.85
$ws.Range("L58").Value = 4983.3335
$ws.Range("M58").Value = -2164.85
$ws.Range("N58").Value = -5389.3335

$ws.Range("H62").Value = 17612.428
$ws.Range("I62").Value = 4077.6
$ws.Range("K62").Value = 4077.6
$ws.Range("M62").Value = -3453.6

$ws.Range("H65").Value = 17612.428
$ws.Range("I65").Value = 4077.6
$ws.Range("K65").Value = 20388
$ws.Range("M65").Value = -17268

$ws.Range("H87").Value = 93000
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 93000
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H97").Value = 41396.184
$ws.Range("I97").Value = 23452.666
$ws.Range("J97").Value = 48125
$ws.Range("K97").Value = 23452.666
$ws.Range("L97").Value = 48125
$ws.Range("M97").Value = -22461.666
$ws.Range("N97").Value = -50107

$ws.Range("H109").Value = 49427
$ws.Range("J109").Value = 48498.25
$ws.Range("L109").Value = 48498.25
$ws.Range("N109").Value = -50578.25

$ws.Range("H132").Value = 74156.82
$ws.Range("I132").Value = 74156.82
$ws.Range("K132").Value = 222470.46
$ws.Range("M132").Value = -219940.46

$ws.Range("H134").Value = 6212.5
$ws.Range("I134").Value = 6255.1
$ws.Range("K134").Value = 18765.3
$ws.Range("M134").Value = -16230.3

$ws.Range("H136").Value = 2709
$ws.Range("I136").Value = 2367.85
$ws.Range("J136").Value = 4983.3335
$ws.Range("K136").Value = 7103.549999999999
$ws.Range("L136").Value = 14950.0005
$ws.Range("M136").Value = -4553.549999999999
$ws.Range("N136").Value = -20050.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 544.5714
$ws.Range("I12").Value = 253.5
$ws.Range("J12").Value = 661
$ws.Range("K12").Value = 760.5
$ws.Range("L12").Value = 1983
$ws.Range("M12").Value = -587.5
$ws.Range("N12").Value = -2329

$ws.Range("H38").Value = 78.1
$ws.Range("I38").Value = 30
$ws.Range("J38").Value = 98.71429
$ws.Range("K38").Value = 90
$ws.Range("L38").Value = 296.14287
$ws.Range("M38").Value = 257
$ws.Range("N38").Value = -990.14287

$ws.Range("H105").Value = 3000
$ws.Range("J105").Value = 3000
$ws.Range("L105").Value = 9000
$ws.Range("N105").Value = -14242

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1827.7354
$ws.Range("I97").Value = 1384.3182
$ws.Range("J97").Value = 2640.6667
$ws.Range("K97").Value = 1384.3182
$ws.Range("L97").Value = 2640.6667
$ws.Range("M97").Value = -888.3181999999999
$ws.Range("N97").Value = -3632.6667

$ws.Range("H122").Value = 3023.7778
$ws.Range("I122").Value = 3087.8572
$ws.Range("K122").Value = 9263.5716
$ws.Range("M122").Value = -6813.571599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 2339.8
$ws.Range("I30").Value = 2339.8
$ws.Range("K30").Value = 2339.8
$ws.Range("M30").Value = -2231.8

$ws.Range("H42").Value = 23875

$ws.Range("H46").Value = 1569.5358
$ws.Range("I46").Value = 772.35
$ws.Range("K46").Value = 772.35
$ws.Range("M46").Value = -584.35

$ws.Range("H49").Value = 23875

$ws.Range("H61").Value = 4152.4707
$ws.Range("I61").Value = 4053.3635
$ws.Range("J61").Value = 4334.1665
$ws.Range("K61").Value = 4053.3635
$ws.Range("L61").Value = 4334.1665
$ws.Range("M61").Value = -3851.3635
$ws.Range("N61").Value = -4738.1665

$ws.Range("H93").Value = 1631.04
$ws.Range("I93").Value = 1284.5883
$ws.Range("J93").Value = 2367.25
$ws.Range("K93").Value = 1284.5883
$ws.Range("L93").Value = 2367.25
$ws.Range("M93").Value = -36.58829999999989
$ws.Range("N93").Value = -4863.25

$ws.Range("H113").Value = 4152.4707
$ws.Range("I113").Value = 4053.3635
$ws.Range("J113").Value = 4334.1665
$ws.Range("K113").Value = 4053.3635
$ws.Range("L113").Value = 4334.1665
$ws.Range("M113").Value = -1883.3635
$ws.Range("N113").Value = -8674.1665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 142862200
$ws.Range("J41").Value = 6639
$ws.Range("L41").Value = 6639
$ws.Range("N41").Value = -7419

$ws.Range("H93").Value = 80666.664

$ws.Range("H109").Value = 84000
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H122").Value = 34519210
$ws.Range("I122").Value = 41709172
$ws.Range("J122").Value = 7379
$ws.Range("K122").Value = 125127516
$ws.Range("L122").Value = 22137
$ws.Range("M122").Value = -125125066
$ws.Range("N122").Value = -27037

$ws.Range("H126").Value = 6650.222
$ws.Range("I126").Value = 8379.23
$ws.Range("K126").Value = 25137.69
$ws.Range("M126").Value = -22667.69

$ws.Range("H132").Value = 1615.8861
$ws.Range("I132").Value = 1448.1493
$ws.Range("K132").Value = 4344.4479
$ws.Range("M132").Value = -1814.4479
